# chore: end of work
# end of work on day 23, 11, 2022
#
# Fill in the missing clock-out data for row 63 (date 2022-11-23 / Alpha EdTech
# studies log): HORA F (clock-out time) and DESCANSO (break time) were blank;
# this was the end-of-day entry that completes them, which also updates the
# ASSUNTO and PRODUÇÃO (subject/output) notes for that day to mention HARD.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estudos")

# HORA F (column C) -> 22:20 and DESCANSO (column E) -> 2:48
$ws.Range("C63").Value = 0.93055555555555547
$ws.Range("E63").Value = 0.11666666666666665

# ASSUNTO (column G) and PRODUÇÃO (column H) updated to reflect the HARD work
$ws.Range("G63").Value = "ESTÁGIO + HARD"
$ws.Range("H63").Value = "Estágio + Atividade voluntária no Alpha EdTech + HARD"
